$wb = $excel.ActiveWorkbook

# Rename existing sheet from "Sheet1" to "Estimates"
$wsEstimates = $wb.Worksheets.Item("Sheet1")
$wsEstimates.Name = "Estimates"

# Add a new sheet "Effort" right after "Estimates"
$wsEffort = $wb.Worksheets.Add([System.Type]::Missing, $wsEstimates)
$wsEffort.Name = "Effort"

# Populate the Effort sheet with header + data row
$wsEffort.Range("A1").Value = "year"
$wsEffort.Range("B1").Value = "n_days"
$wsEffort.Range("A2").Value = 2003
$wsEffort.Range("B2").Value = 1965

# Leave "Effort" as the active/selected sheet, with B3 selected
$wsEffort.Range("B3").Select()
